$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Teacher")

# StudentMaster: ClassId for the demo row was stored as text "1st"; change
# it to the numeric value 1.
$ws.Range("M2").Value = 1

# TeacherAttendancemaster: move the active selection from AE2 to M9.
$ws.Range("M9").Select()
